$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A8').Value = 'Alejandra Hidalgo_20251120_040641'

$ws.Range('B8').NumberFormat = "@"
$ws.Range('B8').Value = ''

$ws.Range('C8').Value = 'Alejandra Hidalgo'

$ws.Range('D8').Value = 21

$ws.Range('E8').Value = 'Female'

$ws.Range('F8').Value = '2025-11-20 04:06:42'

$ws.Range('G8').Value = '{
  "portion": 0.8,
  "diet": 1.0,
  "salt": 0.8,
  "fat": 0.8,
  "natural": 1.0,
  "convenience": 0.6,
  "price": 1.0
}'

$ws.Range('H8').Value = 'Nongshim Neoguri Spicy Seafood'

$ws.Range('I8').NumberFormat = "@"
$ws.Range('I8').Value = '0.569'

$ws.Range('J8').Value = 'Sabor a marisco, umami, picante equilibrado, buena textura, algo salado'

$ws.Range('K8').Value = 'Nissin Chow Mein Teriyaki Beef'

$ws.Range('L8').NumberFormat = "@"
$ws.Range('L8').Value = '0.484'

$ws.Range('M8').Value = 'Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa'

$ws.Range('N8').Value = 'Maruchan Ramen Sabor Pollo'

$ws.Range('O8').NumberFormat = "@"
$ws.Range('O8').Value = '0.467'

$ws.Range('P8').Value = 'Sabor clásico, económico, alto en sodio, no saludable, nostálgico'

$ws.Range('Q8').Value = 'Kraft Macaroni & Cheese Dinner'

$ws.Range('R8').NumberFormat = "@"
$ws.Range('R8').Value = '0.640'

$ws.Range('S8').Value = 'Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato'

$ws.Range('T8').Value = 'Amy’s Macaroni & Cheese (frozen)'

$ws.Range('U8').NumberFormat = "@"
$ws.Range('U8').Value = '0.602'

$ws.Range('V8').Value = 'Queso real, textura casera, sin conservadores, alto en grasa, algo caro'

$ws.Range('W8').Value = 'Annie’s Shells & White Cheddar'

$ws.Range('X8').NumberFormat = "@"
$ws.Range('X8').Value = '0.582'

$ws.Range('Y8').Value = 'Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños'

$ws.Range('Z8').Value = 'Wild Planet Wild Tuna Pasta Salad'

$ws.Range('AA8').NumberFormat = "@"
$ws.Range('AA8').Value = '0.700'

$ws.Range('AB8').Value = 'Sabor fresco, buena proteína, saludable, porción algo pequeña'

$ws.Range('AC8').Value = 'StarKist Chicken Creations (Chicken Salad)'

$ws.Range('AD8').NumberFormat = "@"
$ws.Range('AD8').Value = '0.572'

$ws.Range('AE8').Value = 'Portátil, saludable, fácil, buena textura, sabor suave'

$ws.Range('AF8').Value = 'Jack Link’s Beef Jerky Original'

$ws.Range('AG8').NumberFormat = "@"
$ws.Range('AG8').Value = '0.551'

$ws.Range('AH8').Value = 'Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña'
